# Add new power plants to the "Electricity Source" subscript on the
# BFoCPAbS-electricity sheet (issues #280 and #99).
#
# The sheet lists one electricity-generating technology per row (column A)
# with yearly (2021-2050, columns B:AE) values. We append six new
# technologies as zero-filled rows right after the existing last row (18),
# extending the used range from A1:AE18 to A1:AE24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BFoCPAbS-electricity")

$newTechnologies = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$firstNewRow = 19
$lastCol = 31   # column AE (years 2021-2050 live in B:AE)

$row = $firstNewRow
foreach ($technology in $newTechnologies) {
    $ws.Cells.Item($row, 1).Value() = $technology
    for ($col = 2; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value() = 0
    }
    $row++
}
